$wb = $excel.ActiveWorkbook

# Insert a new worksheet named "total_concentrations" right after
# "input_concentrations" (i.e. before "equilibrium_concentrations"),
# so the sheet order becomes:
#   input_stoich_coefficients, input_k_constants_log10, input_concentrations,
#   total_concentrations, equilibrium_concentrations, PO4_fractions,
#   percent_error, component_names

$before = $wb.Worksheets.Item("equilibrium_concentrations")
$newSheet = $wb.Worksheets.Add($before)
$newSheet.Name = "total_concentrations"

# Populate header row
$newSheet.Range("A1").Value = "H"
$newSheet.Range("B1").Value = "PO4"
$newSheet.Range("C1").Value = "Cu"

# Populate data rows
$newSheet.Range("A2").Value = 0.01
$newSheet.Range("B2").Value = 0.01
$newSheet.Range("C2").Value = 0.01

$newSheet.Range("A3").Value = 0.02
$newSheet.Range("B3").Value = 0.01
$newSheet.Range("C3").Value = 0.01

$newSheet.Range("A4").Value = 0.03
$newSheet.Range("B4").Value = 0.01
$newSheet.Range("C4").Value = 0.01
